$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28 (ALC)
$ws.Range("H28").Value = 443
$ws.Range("I28").Value = 97.333336
$ws.Range("K28").Value = 97.333336
$ws.Range("M28").Value = 387.666664

# Row 33 (ALC)
$ws.Range("H33").Value = 506.625
$ws.Range("I33").Value = 511
$ws.Range("J33").Value = 487.66666
$ws.Range("K33").Value = 511
$ws.Range("L33").Value = 487.66666
$ws.Range("M33").Value = -282
$ws.Range("N33").Value = -945.66666

# Row 96 (ALC)
$ws.Range("H96").Value = 760.86664
$ws.Range("I96").Value = 693.75
$ws.Range("J96").Value = 785.2727
$ws.Range("K96").Value = 2081.25
$ws.Range("L96").Value = 2355.8181
$ws.Range("M96").Value = -708.25
$ws.Range("N96").Value = -5101.8181

# Row 107 (ALC)
$ws.Range("H107").Value = 75481.664
$ws.Range("I107").Value = 100508.336
$ws.Range("K107").Value = 100508.336
$ws.Range("M107").Value = -98588.336

$ws = $wb.Worksheets.Item("ARM")
# Row 18 (ARM)
$ws.Range("H18").Value = 750
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 750
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 750
$ws.Range("N18").Value = -1394
$ws.Range("M18").ClearContents()

# Row 92 (ARM)
$ws.Range("H92").Value = 131499.83
$ws.Range("J92").Value = 131499.83
$ws.Range("L92").Value = 131499.83
$ws.Range("N92").Value = -136491.83

# Row 101 (ARM)
$ws.Range("H101").Value = 9998
$ws.Range("J101").Value = 9998
$ws.Range("L101").Value = 9998
$ws.Range("N101").Value = -16488

# Row 110 (ARM)
$ws.Range("H110").Value = 50000900
$ws.Range("I110").Value = 502.33334
$ws.Range("K110").Value = 502.33334
$ws.Range("M110").Value = 1542.66666

# Row 132 (ARM)
$ws.Range("H132").Value = 2383
$ws.Range("I132").Value = 2400.125
$ws.Range("J132").Value = 2246
$ws.Range("K132").Value = 7200.375
$ws.Range("L132").Value = 6738
$ws.Range("M132").Value = -4670.375
$ws.Range("N132").Value = -11798

$ws = $wb.Worksheets.Item("BSM")
# Row 29 (BSM)
$ws.Range("H29").Value = 1050.75
$ws.Range("I29").Value = 1116.1428
$ws.Range("J29").Value = 999.8889
$ws.Range("K29").Value = 1116.1428
$ws.Range("L29").Value = 999.8889
$ws.Range("M29").Value = -827.1428000000001
$ws.Range("N29").Value = -1577.8889

$ws = $wb.Worksheets.Item("CRP")
# Row 19 (CRP)
$ws.Range("H19").Value = 6000235.5
$ws.Range("I19").Value = 6666900.5
$ws.Range("J19").Value = 250
$ws.Range("K19").Value = 6666900.5
$ws.Range("L19").Value = 250
$ws.Range("M19").Value = -6666730.5
$ws.Range("N19").Value = -590

# Row 22 (CRP)
$ws.Range("H22").Value = 651.0714
$ws.Range("I22").Value = 783.125
$ws.Range("J22").Value = 475
$ws.Range("K22").Value = 783.125
$ws.Range("L22").Value = 475
$ws.Range("M22").Value = -433.125
$ws.Range("N22").Value = -1175

# Row 24 (CRP)
$ws.Range("H24").Value = 6000235.5
$ws.Range("I24").Value = 6666900.5
$ws.Range("J24").Value = 250
$ws.Range("K24").Value = 6666900.5
$ws.Range("L24").Value = 250
$ws.Range("M24").Value = -6666730.5
$ws.Range("N24").Value = -590

# Row 25 (CRP)
$ws.Range("H25").Value = 3333.3333
$ws.Range("I25").Value = 2500
$ws.Range("K25").Value = 2500
$ws.Range("M25").Value = -2326

# Row 31 (CRP)
$ws.Range("H31").Value = 2570.8333
$ws.Range("I31").Value = 1285.1
$ws.Range("J31").Value = 8999.5
$ws.Range("K31").Value = 1285.1
$ws.Range("L31").Value = 8999.5
$ws.Range("M31").Value = -990.0999999999999
$ws.Range("N31").Value = -9589.5

# Row 34 (CRP)
$ws.Range("H34").Value = 2570.8333
$ws.Range("I34").Value = 1285.1
$ws.Range("J34").Value = 8999.5
$ws.Range("K34").Value = 1285.1
$ws.Range("L34").Value = 8999.5
$ws.Range("M34").Value = -1083.1
$ws.Range("N34").Value = -9403.5

# Row 58 (CRP)
$ws.Range("H58").Value = 1000
$ws.Range("I58").Value = 1000
$ws.Range("K58").Value = 1000
$ws.Range("M58").Value = -797

# Row 107 (CRP)
$ws.Range("H107").Value = 1150
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 1300
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 1300
$ws.Range("M107").Value = 920
$ws.Range("N107").Value = -5140

# Row 132 (CRP)
$ws.Range("H132").Value = 1720.75
$ws.Range("I132").Value = 1738
$ws.Range("K132").Value = 5214
$ws.Range("M132").Value = -2684

# Row 134 (CRP)
$ws.Range("H134").Value = 824.75
$ws.Range("I134").Value = 824.75
$ws.Range("K134").Value = 2474.25
$ws.Range("M134").Value = 60.75

# Row 136 (CRP)
$ws.Range("H136").Value = 1000
$ws.Range("I136").Value = 1000
$ws.Range("K136").Value = 3000
$ws.Range("M136").Value = -450

$ws = $wb.Worksheets.Item("CUL")
# Row 12 (CUL)
$ws.Range("H12").Value = 136.53334
$ws.Range("J12").Value = 204.3
$ws.Range("L12").Value = 612.9000000000001
$ws.Range("N12").Value = -958.9000000000001

# Row 139 (CUL)
$ws.Range("H139").Value = 72855.57000000001

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (LTW)
$ws.Range("H22").Value = 1191.7
$ws.Range("I22").Value = 870.44446
$ws.Range("K22").Value = 870.44446
$ws.Range("M22").Value = -575.44446

# Row 27 (LTW)
$ws.Range("H27").Value = 1191.7
$ws.Range("I27").Value = 870.44446
$ws.Range("K27").Value = 870.44446
$ws.Range("M27").Value = -763.44446

# Row 40 (LTW)
$ws.Range("H40").Value = 1171.25
$ws.Range("I40").Value = 1063
$ws.Range("J40").Value = 1496
$ws.Range("K40").Value = 1063
$ws.Range("L40").Value = 1496
$ws.Range("M40").Value = -927
$ws.Range("N40").Value = -1768

# Row 122 (LTW)
$ws.Range("H122").Value = 3493.875
$ws.Range("I122").Value = 3384.7693
$ws.Range("J122").Value = 3966.6667
$ws.Range("K122").Value = 10154.3079
$ws.Range("L122").Value = 11900.0001
$ws.Range("M122").Value = -7704.3079
$ws.Range("N122").Value = -16800.0001

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (WVR)
$ws.Range("H81").Value = 4633.1665
$ws.Range("I81").Value = 4659.8
$ws.Range("J81").Value = 4500
$ws.Range("K81").Value = 9319.6
$ws.Range("L81").Value = 9000
$ws.Range("M81").Value = -8258.6
$ws.Range("N81").Value = -11122

# Row 84 (WVR)
$ws.Range("H84").Value = 4633.1665
$ws.Range("I84").Value = 4659.8
$ws.Range("J84").Value = 4500
$ws.Range("K84").Value = 46598
$ws.Range("L84").Value = 45000
$ws.Range("M84").Value = -41294
$ws.Range("N84").Value = -55608

# Row 122 (WVR)
$ws.Range("H122").Value = 1418.8
$ws.Range("I122").Value = 1032.3334
$ws.Range("J122").Value = 1998.5
$ws.Range("K122").Value = 3097.0002
$ws.Range("L122").Value = 5995.5
$ws.Range("M122").Value = -647.0001999999999
$ws.Range("N122").Value = -10895.5
